$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.451944333333334
$ws.Range("H2").Value = 7.355833000000001
$ws.Range("I2").Value = 0.02621625561007551
$ws.Range("J2").Value = 0.02621625561007551
$ws.Range("Q2").Value = 3.746089542929222
$ws.Range("R2").Value = 33.714805886363
$ws.Range("S2").Value = 0.02621625561007551
$ws.Range("T2").Value = 0.02621625561007551

# Row 3
$ws.Range("I3").Value = 0.02983680056709372
$ws.Range("J3").Value = 0.02983680056709372
$ws.Range("S3").Value = 0.02983680056709372
$ws.Range("T3").Value = 0.02983680056709372

# Row 4
$ws.Range("G4").Value = 24.02567266666667
$ws.Range("H4").Value = 72.07701800000001
$ws.Range("I4").Value = 0.2568831466810099
$ws.Range("J4").Value = 0.25688314668101
$ws.Range("Q4").Value = 36.70651079426644
$ws.Range("R4").Value = 330.358597148398
$ws.Range("S4").Value = 0.2568831466810099
$ws.Range("T4").Value = 0.25688314668101

# Row 5
$ws.Range("G5").Value = 1.586978666666667
$ws.Range("H5").Value = 4.760936
$ws.Range("I5").Value = 0.01696801913790191
$ws.Range("J5").Value = 0.01696801913790192
$ws.Range("Q5").Value = 2.424591825855111
$ws.Range("R5").Value = 21.821326432696
$ws.Range("S5").Value = 0.01696801913790191
$ws.Range("T5").Value = 0.01696801913790192

# Row 6
$ws.Range("G6").Value = 49.62511566666666
$ws.Range("H6").Value = 148.875347
$ws.Range("I6").Value = 0.5305930886400884
$ws.Range("J6").Value = 0.5305930886400885
$ws.Range("Q6").Value = 75.81743367429075
$ws.Range("R6").Value = 682.3569030686169
$ws.Range("S6").Value = 0.5305930886400884
$ws.Range("T6").Value = 0.5305930886400885

# Row 7
$ws.Range("G7").Value = 13.04735633333333
$ws.Range("H7").Value = 39.142069
$ws.Range("I7").Value = 0.1395026893638304
$ws.Range("J7").Value = 0.1395026893638304
$ws.Range("Q7").Value = 19.93379884637322
$ws.Range("R7").Value = 179.404189617359
$ws.Range("S7").Value = 0.1395026893638304
$ws.Range("T7").Value = 0.1395026893638304
